$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1812.04
$ws.Range("I41").Value = 2262
$ws.Range("J41").Value = 1012.1111
$ws.Range("K41").Value = 2262
$ws.Range("L41").Value = 1012.1111
$ws.Range("M41").Value = -1822
$ws.Range("N41").Value = -1892.1111
$ws.Range("H61").Value = 180952380
$ws.Range("I61").Value = 28571428
$ws.Range("K61").Value = 85714284
$ws.Range("M61").Value = -85714112
$ws.Range("H62").Value = 22226428
$ws.Range("I62").Value = 40004370
$ws.Range("J62").Value = 3999.75
$ws.Range("K62").Value = 40004370
$ws.Range("L62").Value = 3999.75
$ws.Range("M62").Value = -40003746
$ws.Range("N62").Value = -5247.75
$ws.Range("H65").Value = 22226428
$ws.Range("I65").Value = 40004370
$ws.Range("J65").Value = 3999.75
$ws.Range("K65").Value = 200021850
$ws.Range("L65").Value = 19998.75
$ws.Range("M65").Value = -200018730
$ws.Range("N65").Value = -26238.75
$ws.Range("H86").Value = 4334.3335
$ws.Range("I86").Value = 4001.5
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4001.5
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2878.5
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4334.3335
$ws.Range("I89").Value = 4001.5
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 20007.5
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -14391.5
$ws.Range("N89").Value = -36232
$ws.Range("H98").Value = 1281.2941
$ws.Range("J98").Value = 1899.3334
$ws.Range("L98").Value = 1899.3334
$ws.Range("N98").Value = -4895.3334
$ws.Range("H122").Value = 1281.2941
$ws.Range("J122").Value = 1899.3334
$ws.Range("L122").Value = 5698.0002
$ws.Range("N122").Value = -10598.0002
$ws.Range("H132").Value = 64106.125
$ws.Range("I132").Value = 39949.883
$ws.Range("J132").Value = 168783.17
$ws.Range("K132").Value = 119849.649
$ws.Range("L132").Value = 506349.51
$ws.Range("M132").Value = -117319.649
$ws.Range("N132").Value = -511409.51
$ws.Range("H137").Value = 575939.0600000001
$ws.Range("I137").Value = 21034.924
$ws.Range("K137").Value = 63104.772
$ws.Range("M137").Value = -60554.772
$ws.Range("H138").Value = 3274.383
$ws.Range("I138").Value = 1239.4286
$ws.Range("J138").Value = 4137.697
$ws.Range("K138").Value = 3718.2858
$ws.Range("L138").Value = 12413.091
$ws.Range("M138").Value = 1421.7142
$ws.Range("N138").Value = -22693.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 125003960
$ws.Range("J2").Value = 2500
$ws.Range("L2").Value = 2500
$ws.Range("N2").Value = -2726
$ws.Range("H32").Value = 37130.73
$ws.Range("I32").Value = 27677.176
$ws.Range("K32").Value = 27677.176
$ws.Range("M32").Value = -27390.176
$ws.Range("H61").Value = 3512.6765
$ws.Range("I61").Value = 3407
$ws.Range("K61").Value = 3407
$ws.Range("M61").Value = -3195
$ws.Range("H88").Value = 2615.9443
$ws.Range("I88").Value = 2613
$ws.Range("K88").Value = 2613
$ws.Range("M88").Value = -2207
$ws.Range("H91").Value = 2615.9443
$ws.Range("I91").Value = 2613
$ws.Range("K91").Value = 2613
$ws.Range("M91").Value = -1209
$ws.Range("H116").Value = 125003960
$ws.Range("J116").Value = 2500
$ws.Range("L116").Value = 2500
$ws.Range("N116").Value = -7088
$ws.Range("H122").Value = 1039.25
$ws.Range("I122").Value = 887.3333
$ws.Range("K122").Value = 2661.9999
$ws.Range("M122").Value = -211.9998999999998
$ws.Range("H132").Value = 3013.6
$ws.Range("I132").Value = 2421.389
$ws.Range("J132").Value = 4536.4287
$ws.Range("K132").Value = 7264.167
$ws.Range("L132").Value = 13609.2861
$ws.Range("M132").Value = -4734.167
$ws.Range("N132").Value = -18669.2861
$ws.Range("H136").Value = 3512.6765
$ws.Range("I136").Value = 3407
$ws.Range("K136").Value = 10221
$ws.Range("M136").Value = -7671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 125003960
$ws.Range("J3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("N3").Value = -2728
$ws.Range("H105").Value = 7926.625
$ws.Range("J105").Value = 14337
$ws.Range("L105").Value = 14337
$ws.Range("N105").Value = -17831
$ws.Range("H107").Value = 1274.6
$ws.Range("I107").Value = 1222.6
$ws.Range("J107").Value = 1378.6
$ws.Range("K107").Value = 1222.6
$ws.Range("L107").Value = 1378.6
$ws.Range("M107").Value = 697.4000000000001
$ws.Range("N107").Value = -5218.6
$ws.Range("H134").Value = 1992.6552
$ws.Range("I134").Value = 1818.8846
$ws.Range("K134").Value = 5456.6538
$ws.Range("M134").Value = -2921.6538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7848912
$ws.Range("I31").Value = 4221886.5
$ws.Range("J31").Value = 13893955
$ws.Range("K31").Value = 4221886.5
$ws.Range("L31").Value = 13893955
$ws.Range("M31").Value = -4221591.5
$ws.Range("N31").Value = -13894545
$ws.Range("H34").Value = 7848912
$ws.Range("I34").Value = 4221886.5
$ws.Range("J34").Value = 13893955
$ws.Range("K34").Value = 4221886.5
$ws.Range("L34").Value = 13893955
$ws.Range("M34").Value = -4221684.5
$ws.Range("N34").Value = -13894359
$ws.Range("H105").Value = 900
$ws.Range("I105").Value = 900
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 900
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 847
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 812.5172
$ws.Range("I107").Value = 825
$ws.Range("K107").Value = 825
$ws.Range("M107").Value = 1095
$ws.Range("H122").Value = 2397.2778
$ws.Range("I122").Value = 2338.3572
$ws.Range("K122").Value = 7015.071599999999
$ws.Range("M122").Value = -4565.071599999999
$ws.Range("H132").Value = 1610.2572
$ws.Range("I132").Value = 1613.5
$ws.Range("K132").Value = 4840.5
$ws.Range("M132").Value = -2310.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 3131
$ws.Range("I8").Value = 3131
$ws.Range("K8").Value = 9393
$ws.Range("M8").Value = -9254
$ws.Range("H121").Value = 3077.2222
$ws.Range("I121").Value = 2036.9546
$ws.Range("J121").Value = 7654.4
$ws.Range("K121").Value = 6110.8638
$ws.Range("L121").Value = 22963.2
$ws.Range("M121").Value = -4800.8638
$ws.Range("N121").Value = -25583.2
$ws.Range("H129").Value = 5824431
$ws.Range("I129").Value = 9000393
$ws.Range("J129").Value = 1833.3334
$ws.Range("K129").Value = 27001179
$ws.Range("L129").Value = 5500.0002
$ws.Range("M129").Value = -26996179
$ws.Range("N129").Value = -15500.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 18000
$ws.Range("I46").Value = 35000
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 35000
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -34844
$ws.Range("N46").Value = -1312
$ws.Range("H80").Value = 166668140
$ws.Range("J80").Value = 149
$ws.Range("L80").Value = 149
$ws.Range("N80").Value = -2145
$ws.Range("H83").Value = 166668140
$ws.Range("J83").Value = 149
$ws.Range("L83").Value = 745
$ws.Range("N83").Value = -10729
$ws.Range("H97").Value = 1575.9656
$ws.Range("I97").Value = 1507.25
$ws.Range("K97").Value = 1507.25
$ws.Range("M97").Value = -1011.25
$ws.Range("H132").Value = 44385.31
$ws.Range("J132").Value = 1415.5
$ws.Range("L132").Value = 4246.5
$ws.Range("N132").Value = -9306.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5020.7856
$ws.Range("I122").Value = 4537.25
$ws.Range("J122").Value = 5665.5
$ws.Range("K122").Value = 13611.75
$ws.Range("L122").Value = 16996.5
$ws.Range("M122").Value = -11161.75
$ws.Range("N122").Value = -21896.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 807.1070999999999
$ws.Range("I107").Value = 620.4545000000001
$ws.Range("J107").Value = 1491.5
$ws.Range("K107").Value = 1861.3635
$ws.Range("L107").Value = 4474.5
$ws.Range("M107").Value = 58.63649999999984
$ws.Range("N107").Value = -8314.5
$ws.Range("H122").Value = 3185.6843
$ws.Range("I122").Value = 3189.25
$ws.Range("K122").Value = 9567.75
$ws.Range("M122").Value = -7117.75
$ws.Range("H132").Value = 3410.3667
$ws.Range("I132").Value = 3592.6155
$ws.Range("K132").Value = 10777.8465
$ws.Range("M132").Value = -8247.8465
$ws.Range("H136").Value = 1976.6364
$ws.Range("I136").Value = 1730.375
$ws.Range("K136").Value = 5191.125
$ws.Range("M136").Value = -2641.125
